# Natmi following Dr Hou advice
# Adds the "ECs" sending/target cluster to the Bmp4 -> Rgmb NATMI LR-pair sheet,
# expanding the existing FAPs/sCs x FAPs/sCs grid (rows 2-7) into a full
# ECs/FAPs/sCs x ECs/FAPs/sCs grid (rows 2-10) with refreshed NATMI statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs (Bmp4/Rgmb)
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Bmp4"
$ws.Range("C2").Value = "Rgmb"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 5.269639666666666
$ws.Range("H2").Value = 15.808919
$ws.Range("I2").Value = 0.09922110188645328
$ws.Range("J2").Value = 0.09922110188645328
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.367553
$ws.Range("N2").Value = 13.102659
$ws.Range("O2").Value = 0.05657159077620311
$ws.Range("P2").Value = 0.05657159077620311
$ws.Range("Q2").Value = 23.015430535069
$ws.Range("R2").Value = 207.138874815621
$ws.Range("S2").Value = 0.005613095572284389
$ws.Range("T2").Value = 0.005613095572284389

# Row 3: ECs -> FAPs (Bmp4/Rgmb)
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Bmp4"
$ws.Range("C3").Value = "Rgmb"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 5.269639666666666
$ws.Range("H3").Value = 15.808919
$ws.Range("I3").Value = 0.09922110188645328
$ws.Range("J3").Value = 0.09922110188645328
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 53.45120266666667
$ws.Range("N3").Value = 160.353608
$ws.Range("O3").Value = 0.6923372340884159
$ws.Range("P3").Value = 0.6923372340884159
$ws.Range("Q3").Value = 281.6685778033058
$ws.Range("R3").Value = 2535.017200229752
$ws.Range("S3").Value = 0.06869446324327197
$ws.Range("T3").Value = 0.06869446324327197

# Row 4: ECs -> sCs (Bmp4/Rgmb)
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Bmp4"
$ws.Range("C4").Value = "Rgmb"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 5.269639666666666
$ws.Range("H4").Value = 15.808919
$ws.Range("I4").Value = 0.09922110188645328
$ws.Range("J4").Value = 0.09922110188645328
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 19.38524266666667
$ws.Range("N4").Value = 58.155728
$ws.Range("O4").Value = 0.251091175135381
$ws.Range("P4").Value = 0.251091175135381
$ws.Range("Q4").Value = 102.1532437042258
$ws.Range("R4").Value = 919.379193338032
$ws.Range("S4").Value = 0.02491354307089693
$ws.Range("T4").Value = 0.02491354307089693

# Row 5: FAPs -> ECs (Bmp4/Rgmb)
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Bmp4"
$ws.Range("C5").Value = "Rgmb"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 39.012863
$ws.Range("H5").Value = 117.038589
$ws.Range("I5").Value = 0.7345662131494083
$ws.Range("J5").Value = 0.7345662131494083
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.367553
$ws.Range("N5").Value = 13.102659
$ws.Range("O5").Value = 0.05657159077620311
$ws.Range("P5").Value = 0.05657159077620311
$ws.Range("Q5").Value = 170.390746834239
$ws.Range("R5").Value = 1533.516721508151
$ws.Range("S5").Value = 0.04155557920831351
$ws.Range("T5").Value = 0.04155557920831351

# Row 6: FAPs -> FAPs (Bmp4/Rgmb)
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Bmp4"
$ws.Range("C6").Value = "Rgmb"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 39.012863
$ws.Range("H6").Value = 117.038589
$ws.Range("I6").Value = 0.7345662131494083
$ws.Range("J6").Value = 0.7345662131494083
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 53.45120266666667
$ws.Range("N6").Value = 160.353608
$ws.Range("O6").Value = 0.6923372340884159
$ws.Range("P6").Value = 0.6923372340884159
$ws.Range("Q6").Value = 2085.284446819901
$ws.Range("R6").Value = 18767.56002137911
$ws.Range("S6").Value = 0.5085675402666631
$ws.Range("T6").Value = 0.5085675402666631

# Row 7: FAPs -> sCs (Bmp4/Rgmb)
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Bmp4"
$ws.Range("C7").Value = "Rgmb"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 39.012863
$ws.Range("H7").Value = 117.038589
$ws.Range("I7").Value = 0.7345662131494083
$ws.Range("J7").Value = 0.7345662131494083
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 19.38524266666667
$ws.Range("N7").Value = 58.155728
$ws.Range("O7").Value = 0.251091175135381
$ws.Range("P7").Value = 0.251091175135381
$ws.Range("Q7").Value = 756.2738163764213
$ws.Range("R7").Value = 6806.464347387792
$ws.Range("S7").Value = 0.1844430936744317
$ws.Range("T7").Value = 0.1844430936744317

# Row 8: sCs -> ECs (Bmp4/Rgmb)
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Bmp4"
$ws.Range("C8").Value = "Rgmb"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 8.827567333333334
$ws.Range("H8").Value = 26.482702
$ws.Range("I8").Value = 0.1662126849641383
$ws.Range("J8").Value = 0.1662126849641383
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 4.367553
$ws.Range("N8").Value = 13.102659
$ws.Range("O8").Value = 0.05657159077620311
$ws.Range("P8").Value = 0.05657159077620311
$ws.Range("Q8").Value = 38.55486818940201
$ws.Range("R8").Value = 346.993813704618
$ws.Range("S8").Value = 0.0094029159956052
$ws.Range("T8").Value = 0.0094029159956052

# Row 9: sCs -> FAPs (Bmp4/Rgmb)
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Bmp4"
$ws.Range("C9").Value = "Rgmb"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 8.827567333333334
$ws.Range("H9").Value = 26.482702
$ws.Range("I9").Value = 0.1662126849641383
$ws.Range("J9").Value = 0.1662126849641383
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 53.45120266666667
$ws.Range("N9").Value = 160.353608
$ws.Range("O9").Value = 0.6923372340884159
$ws.Range("P9").Value = 0.6923372340884159
$ws.Range("Q9").Value = 471.8440905876463
$ws.Range("R9").Value = 4246.596815288817
$ws.Range("S9").Value = 0.1150752305784807
$ws.Range("T9").Value = 0.1150752305784807

# Row 10: sCs -> sCs (Bmp4/Rgmb)
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Bmp4"
$ws.Range("C10").Value = "Rgmb"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 8.827567333333334
$ws.Range("H10").Value = 26.482702
$ws.Range("I10").Value = 0.1662126849641383
$ws.Range("J10").Value = 0.1662126849641383
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 19.38524266666667
$ws.Range("N10").Value = 58.155728
$ws.Range("O10").Value = 0.251091175135381
$ws.Range("P10").Value = 0.251091175135381
$ws.Range("Q10").Value = 171.1245349130062
$ws.Range("R10").Value = 1540.120814217056
$ws.Range("S10").Value = 0.04173453839005237
$ws.Range("T10").Value = 0.04173453839005237
